# respostasalunos.xlsx -- "Add files via upload"
#
# The OOXML diff shows a lot of "noise" that Excel regenerates on every
# save (fileVersion/rupBuild, xr:revisionPtr GUID, and a full shared-string
# table reshuffle that simply drops three now-unused strings from the
# front of the table and appends three new ones at the back -- every data
# row's <v> shared-string index shifts by exactly 3, but the *text* each
# row points at is unchanged). The only real, intentional edit is:
#
#   B1: "Período de Interesse"                                            -> "interesse"
#   C1: "Possui computador em casa para realização das aulas práticas?"   -> "computador"
#   D1: "Qual sua motivação pra realização do curso?"                     -> "motivacao"
#
# plus the sheet being left scrolled further down with A67 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text changes (B1:D1) -------------------------------------------------
$ws.Range("B1").Value = "interesse"
$ws.Range("C1").Value = "computador"
$ws.Range("D1").Value = "motivacao"

# --- Selection / scroll state -----------------------------------------------------
# Before: frozen pane topLeftCell A2, selection A4
# After:  frozen pane topLeftCell A52, selection A67
[void]$ws.Range("A67").Select()
$excel.ActiveWindow.ScrollRow = 52
